# Add a new "2022-Q1" worksheet (placed between "2021-Q4" and "总计"),
# populate it with the Q1-2022 fund-holding data, and update the "总计"
# (totals) sheet with a new row summarising the 2022-Q1 period.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Header row text
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Match the header formatting (bold, centered, bordered) already used on
# the "2021-Q4" sheet by copying it over.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Data rows. Numeric-looking values in columns B, D, E, F, G are stored
# as text in the source workbook (fund codes keep leading zeros, and the
# numeric-looking measures are text too), so they are entered with a
# leading apostrophe to force a text interpretation instead of being
# auto-converted to numbers.
$q1Data = @(
    @(0, "013393", "信达澳银价值精选混合A",       "3.61",  "81.31", "5.00", "0.1805", 1),
    @(1, "003655", "信达澳银新财富灵活配置混合",   "11.86", "25.86", "0.78", "0.0925", 1),
    @(2, "003456", "信达澳银新目标灵活配置混合",   "1.05",  "86.04", "1.76", "0.0185", 9),
    @(3, "013394", "信达澳银价值精选混合C",       "0.37",  "81.31", "5.00", "0.0185", 1),
    @(4, "012005", "信达澳银恒盛混合A",           "1.87",  "31.90", "0.81", "0.0151", 6),
    @(5, "012006", "信达澳银恒盛混合C",           "0.31",  "31.90", "0.81", "0.0025", 6)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1.Cells.Item($r, 3).Value = "'" + $row[2]
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Column A on the source sheet uses the bold/centered/bordered style too;
# replicate that down the new data rows.
$q4.Range("A2").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row for 2022-Q1
#    above the existing 2021-Q4 totals row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "'2022-Q1"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.33

$total.Cells.Item(3, 1).Value = 1

# Restore the A-column style (bold/centered/bordered) on the new row to
# match the row below it.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
